# Adds "Rigid 10K" print-settings info to the Blad1 (printers/cutters) sheet.
#
# Structurally this mirrors what Excel does for "Insert Column": three new
# blank columns are inserted before column F (old F:L shift right to I:O),
# and a fourth new column is inserted a little further along (pushing the
# old K:L -- now N:O -- one more step right, to P:Q). The newly freed
# columns (F, G, H and the later N) are then used to hold a small
# "rigid 10K V1" curing-settings box plus a "materials" label.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Blad1")

# --- 1. Insert the new columns -------------------------------------------
# Three blank columns before old column F (old F..J -> new I..M, etc.)
$ws.Range("F1:H1").EntireColumn.Insert()
# One more blank column further right (old K -> new N is pushed to O, etc.)
$ws.Range("N1").EntireColumn.Insert()

# --- 2. New "materials" label ---------------------------------------------
$ws.Range("A13").Value = "materials"
$ws.Range("G13").Value = "materials"

# --- 3. New "rigid 10K V1" curing-settings mini table ----------------------
$ws.Range("H16").Value = "rigid 10K V1"
$ws.Range("I16").Value = "curing time"
$ws.Range("J16").Value = "60 min"
$ws.Range("K16").Value = "~"

$ws.Range("I17").Value = "curing temp"
$ws.Range("J17").Value = "70 °C"
$ws.Range("K17").Value = "~"

$ws.Range("I18").Value = "note"
$ws.Range("J18").Value = "Optionally, post-cure and then heat the printed part at 125 ºC for 90 minutes for a higher heat deflection temperature."
$ws.Range("K18").Value = "~"

# --- 4. Light formatting so the new box reads like the rest of the sheet ---
$box = $ws.Range("H16:K18")
$box.Borders.Item(7).LineStyle = 1   # xlEdgeLeft
$box.Borders.Item(8).LineStyle = 1   # xlEdgeTop
$box.Borders.Item(9).LineStyle = 1   # xlEdgeBottom
$box.Borders.Item(10).LineStyle = 1  # xlEdgeRight
$box.Borders.Item(11).LineStyle = 1  # xlInsideVertical
$box.Borders.Item(12).LineStyle = 1  # xlInsideHorizontal

$ws.Range("J18").WrapText = $true
$ws.Range("H16").Font.Bold = $false

# --- 5. Selection / view bookkeeping (cosmetic, matches the saved view) ---
[void]$ws.Range("K20").Select()
